$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price column so numeric-looking strings
# (trailing zeros, thousand-dot separators, etc.) are kept verbatim.
$priceCells = "D2","D3","D5","D6","D8","D10","D11","D12","D14","D15","D16","D17","D18","D20","D21","D23","D24","D25","D26","D27","D28","D29","D30","D32","D35","D36","D37","D38","D40","D42","D43","D44","D45","D47","D49","D50","D51"
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "66.041.02"
$ws.Range("E2").Value = "  +2.62%  "

# Row 3
$ws.Range("D3").Value = "3.177.73"
$ws.Range("E3").Value = "  +1.78%  "

# Row 4
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").Value = "594.52"
$ws.Range("E5").Value = "  +4.71%  "

# Row 6
$ws.Range("D6").Value = "152.98"
$ws.Range("E6").Value = "  +4.00%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "3.168.98"
$ws.Range("E8").Value = "  +1.47%  "

# Row 9
$ws.Range("E9").Value = "  +4.24%  "

# Row 10
$ws.Range("D10").Value = "0.158"
$ws.Range("E10").Value = "  +1.86%  "

# Row 11
$ws.Range("D11").Value = "6.00"
$ws.Range("E11").Value = "  -0.37%  "

# Row 12
$ws.Range("D12").Value = "0.512"
$ws.Range("E12").Value = "  +4.73%  "

# Row 13
$ws.Range("E13").Value = "  +2.98%  "

# Row 14
$ws.Range("D14").Value = "38.83"
$ws.Range("E14").Value = "  +6.46%  "

# Row 15
$ws.Range("D15").Value = "3.700.81"
$ws.Range("E15").Value = "  +1.86%  "

# Row 16
$ws.Range("D16").Value = "66.042.91"
$ws.Range("E16").Value = "  +2.19%  "

# Row 17
$ws.Range("D17").Value = "7.41"
$ws.Range("E17").Value = "  +6.09%  "

# Row 18
$ws.Range("D18").Value = "3.182.94"
$ws.Range("E18").Value = "  +1.46%  "

# Row 19
$ws.Range("E19").Value = "  +1.15%  "

# Row 20
$ws.Range("D20").Value = "507.63"
$ws.Range("E20").Value = "  +2.39%  "

# Row 21
$ws.Range("D21").Value = "15.24"
$ws.Range("E21").Value = "  +4.50%  "

# Row 22
$ws.Range("E22").Value = "  +4.19%  "

# Row 23
$ws.Range("D23").Value = "7.98"
$ws.Range("E23").Value = "  +5.28%  "

# Row 24
$ws.Range("D24").Value = "15.01"
$ws.Range("E24").Value = "  +0.55%  "

# Row 25
$ws.Range("D25").Value = "84.62"
$ws.Range("E25").Value = "  +1.66%  "

# Row 26
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  +0.13%  "

# Row 27
$ws.Range("D27").Value = "9.16"
$ws.Range("E27").Value = "  +5.08%  "

# Row 28
$ws.Range("D28").Value = "2.99"
$ws.Range("E28").Value = "  +4.75%  "

# Row 29
$ws.Range("D29").Value = "2.27"
$ws.Range("E29").Value = "  +6.61%  "

# Row 30
$ws.Range("D30").Value = "7.00"
$ws.Range("E30").Value = "  +15.27%  "

# Row 31
$ws.Range("E31").Value = "  +5.07%  "

# Row 32
$ws.Range("D32").Value = "28.05"
$ws.Range("E32").Value = "  +3.28%  "

# Row 33
$ws.Range("E33").Value = "  +3.46%  "

# Row 34
$ws.Range("E34").Value = "  -0.21%  "

# Row 35
$ws.Range("D35").Value = "6.46"
$ws.Range("E35").Value = "  +1.57%  "

# Row 36
$ws.Range("D36").Value = "54.68"
$ws.Range("E36").Value = "  +0.95%  "

# Row 37
$ws.Range("D37").Value = "487.26"
$ws.Range("E37").Value = "  +6.07%  "

# Row 38
$ws.Range("D38").Value = "0.0892"
$ws.Range("E38").Value = "  +0.83%  "

# Row 39
$ws.Range("E39").Value = "  +1.92%  "

# Row 40
$ws.Range("D40").Value = "8.84"
$ws.Range("E40").Value = "  +3.90%  "

# Row 41
$ws.Range("E41").Value = "  +5.00%  "

# Row 42
$ws.Range("D42").Value = "0.297"
$ws.Range("E42").Value = "  +6.89%  "

# Row 43
$ws.Range("D43").Value = "2.81"
$ws.Range("E43").Value = "  -3.46%  "

# Row 44 (price contains U+2083 SUBSCRIPT THREE, written here as a literal
# UTF-8 character rather than built via [char]/+ concatenation)
$ws.Range("D44").Value = "0.0₃0653"
$ws.Range("E44").Value = "  +15.65%  "

# Row 45
$ws.Range("D45").Value = "2.899.00"
$ws.Range("E45").Value = "  -3.59%  "

# Row 46
$ws.Range("E46").Value = "  +0.93%  "

# Row 47
$ws.Range("D47").Value = "28.25"
$ws.Range("E47").Value = "  +1.32%  "

# Row 48
$ws.Range("E48").Value = "  +0.02%  "

# Row 49 - swapped with row 50 (CoreDAO <-> Stellar) with new values
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "0.116"
$ws.Range("E49").Value = "  +2.95%  "

# Row 50
$ws.Range("B50").Value = "CoreDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D50").Value = "2.65"
$ws.Range("E50").Value = "  +11.73%  "

# Row 51
$ws.Range("D51").Value = "2.30"
$ws.Range("E51").Value = "  +4.82%  "
